$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Dependencies here" threaded comment from A9 to A14 ---
$commentText = $ws.Range("A9").Comment.Text()
$ws.Range("A9").Comment.Delete()

# --- Write the new version labels first (this controls the order new ---
# --- strings are interned into the shared-strings table) ---
$ws.Range("F4").Value = "V1.03.1"
$ws.Range("F5").Value = "V1.03.2"
$ws.Range("A4").Value = "V1.01.1"
$ws.Range("A3").Value = "V1.02.1"
$ws.Range("A5").Value = "V1.01.2"
$ws.Range("A6").Value = "V1.01.3"
$ws.Range("A7").Value = "V1.01.5"
$ws.Range("A10").Value = "V1.03.6"

# --- Fill in the rest of the numeric data ---

# Row 4: rest of new data (A4/F4 already set above); K4/L4/N4/O4 already existed
$ws.Range("B4").Value = 9.4274299999999993
$ws.Range("C4").Value = 0.11418
$ws.Range("D4").Value = 0.53954999999999997
$ws.Range("G4").Value = 22.260269999999998
$ws.Range("H4").Value = 0.00599
$ws.Range("I4").Value = 0.79791999999999996

# Row 5: rest of new row (A5/F5 already set above)
$ws.Range("B5").Value = 13.53234
$ws.Range("C5").Value = 0.04447
$ws.Range("D5").Value = 0.56920999999999999
$ws.Range("G5").Value = 20.749400000000001
$ws.Range("H5").Value = 0.00851
$ws.Range("I5").Value = 0.82787999999999995

# Row 6: rest of new row (A6 already set above)
$ws.Range("B6").Value = 13.17601
$ws.Range("C6").Value = 0.04825
$ws.Range("D6").Value = 0.53405000000000002

# Row 7: rest of new row (A7 already set above)
$ws.Range("B7").Value = 13.1869
$ws.Range("C7").Value = 0.04811
$ws.Range("D7").Value = 0.40762999999999999

# Row 8: new row
$ws.Range("A8").Value = "V1.03.1"
$ws.Range("B8").Value = 18.560929999999999
$ws.Range("C8").Value = 0.01395
$ws.Range("D8").Value = 0.73333999999999999

# Row 9: label now "V1.03.2" (previously held "Dependencies here" / the comment anchor)
$ws.Range("A9").Value = "V1.03.2"

# Row 10: rest of new row (A10 already set above)
$ws.Range("B10").Value = 18.985990000000001
$ws.Range("C10").Value = 0.01265
$ws.Range("D10").Value = 0.70992999999999995

# Row 14: "Dependencies here" label, now carries the relocated comment
$ws.Range("A14").Value = "Dependencies here"
$ws.Range("A14").AddCommentThreaded($commentText) | Out-Null

# --- Update the active selection to match the new view ---
$ws.Range("C11").Select() | Out-Null
